$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 79, shifting existing rows 79..178 down to 80..179.
$ws.Rows(79).Insert()

# Populate the newly inserted row 79 with the new record.
# Columns that stay the same as the (now shifted-down) row 80 / original row 79:
$ws.Range("A79").Value = 7
$ws.Range("B79").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C79").Value = "Ñuble"
$ws.Range("D79").Value = 44778
$ws.Range("E79").Value = 16
$ws.Range("F79").Value = "Fruta"
$ws.Range("G79").Value = 100102
$ws.Range("H79").Value = "Cítricos"
$ws.Range("I79").Value = 100102004
$ws.Range("J79").Value = "Mandarina"
$ws.Range("K79").Value = "Clemenuless"
$ws.Range("L79").Value = "Primera"
$ws.Range("M79").Value = 160
$ws.Range("N79").Value = 8500
$ws.Range("O79").Value = 9000
$ws.Range("P79").Value = 8750
$ws.Range("Q79").Value = "$/caja 18 kilos"
$ws.Range("R79").Value = "Región de O'Higgins"
$ws.Range("S79").Value = 486
$ws.Range("T79").Value = 18
